$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): rename "Ngày tham gia" -> "Giới tính", add "Ngày sinh" and "Địa chỉ" ---
$ws.Range("E1").Value = "Giới tính"

$ws.Range("F1").Value = "Ngày sinh"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $False

$ws.Range("G1").Value = "Địa chỉ"
$ws.Range("E1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = $False

# --- Remove old hyperlinks so we can recreate them with fresh targets ---
$ws.Hyperlinks.Delete()

# --- Row 2: Nguyễn Đình Hùng ---
$ws.Range("A2").Value = 12520011
$ws.Range("B2").Value = "Nguyễn Đình Hùng"
$ws.Range("C2").Value = "nguyendinhhungzz@gmail.com"
$ws.Range("D2").Value = 2385785342
$ws.Range("E2").Value = "Nam"
$ws.Range("F2").Value = 37535
$ws.Range("F2").NumberFormat = "mm-dd-yy"
$ws.Range("G2").Value = "Hưng Yên"

# --- Row 3: Phạm Thanh Long ---
$ws.Range("A3").Value = 12520012
$ws.Range("B3").Value = "Phạm Thanh Long"
$ws.Range("C3").Value = "phamthanhlongê@gmail.com"
$ws.Range("D3").Value = 1385785342
$ws.Range("E3").Value = "Nam"
$ws.Range("F3").Value = 37562
$ws.Range("F3").NumberFormat = "mm-dd-yy"
$ws.Range("G3").Value = "Hưng Yên"

# --- Row 4: Bùi Xuân Hoàng ---
$ws.Range("A4").Value = 12520013
$ws.Range("B4").Value = "Bùi Xuân Hoàng"
$ws.Range("C4").Value = "buixuanhoangyy@gmail.com"
$ws.Range("D4").Value = 5385785342
$ws.Range("E4").Value = "Nữ"
$ws.Range("F4").Value = 37351
$ws.Range("F4").NumberFormat = "mm-dd-yy"
$ws.Range("G4").Value = "Hà Nội"

# --- Re-add mail hyperlinks on the Email column ---
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:nguyendinhhungzz@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:phamthanhlongê@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:buixuanhoangyy@gmail.com")

# --- Column widths for the new/changed columns ---
$ws.Range("E1").ColumnWidth = 8.43
$ws.Range("F1").ColumnWidth = 9.25
$ws.Range("G1").ColumnWidth = 8.25

# --- Selection state matches the saved file ---
[void]$ws.Range("D8").Select()
